$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Model R^2
$ws.Range("B7").Value = 0.218
$ws.Range("C7").Value = 0.676
$ws.Range("D7").Value = 0.468
$ws.Range("E7").Value = 0.628
$ws.Range("F7").Value = 0.462
$ws.Range("G7").Value = 0.515

# Row 8 - Model Adj R^2
$ws.Range("B8").Value = -16.055
$ws.Range("C8").Value = 0.453
$ws.Range("D8").Value = -1.659
$ws.Range("E8").Value = 0.574
$ws.Range("F8").Value = 0.441
$ws.Range("G8").Value = 0.502

# Row 9 - Model RMSE
$ws.Range("B9").Value = 2.633
$ws.Range("C9").Value = 2.012
$ws.Range("D9").Value = 2.292
$ws.Range("E9").Value = 2.31
$ws.Range("F9").Value = 2.426
$ws.Range("G9").Value = 2.381

# Row 10 - Model HH
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 18
$ws.Range("G10").Value = 67

# Row 11 - Delta R^2
$ws.Range("B11").Value = 0.886
$ws.Range("C11").Value = 0.363
$ws.Range("D11").Value = 0.525
$ws.Range("E11").Value = 0.63
$ws.Range("F11").Value = 0.495
$ws.Range("G11").Value = 0.526

# Row 12 - Delta Adj R^2
$ws.Range("B12").Value = 19.325
$ws.Range("C12").Value = 0.613
$ws.Range("D12").Value = 2.626
$ws.Range("E12").Value = 0.722
$ws.Range("F12").Value = 0.514
$ws.Range("G12").Value = 0.539

# Row 13 - Delta RMSE
$ws.Range("B13").Value = -1.213
$ws.Range("C13").Value = -0.918
$ws.Range("D13").Value = -0.939
$ws.Range("E13").Value = -1.48
$ws.Range("F13").Value = -0.936
$ws.Range("G13").Value = -1.056

# Row 14 - Delta HH
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = -5
$ws.Range("G14").Value = -7
